$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "MVP overzicht": update a handful of "Issue(s)" cells in column E
# ---------------------------------------------------------------------------
$wsOverzicht = $wb.Worksheets.Item("MVP overzicht")

# Row 21 (item #8): append issue 480 to the list
$wsOverzicht.Range("E21").Value = "637 (M)638, 639, 481(M), 532, 531, 480"

# Row 24 (item #11): previously empty, now holds a new issue number
$wsOverzicht.Range("E24").Value = 474

# Row 37 (item 'Als medewerker wil ik een Babs kunnen toevoegen/bewerken/verwijderen'):
# was a bare number (432), now a text list including the new issue 477
$wsOverzicht.Range("E37").Value = "432, 477"

# Row 39 (item #23): append issue 492 to the list
$wsOverzicht.Range("E39").Value = "561(L), 492"

# Row 58 (item #38 'Koppeling met openzaak'): was a bare number (476), now a
# text list including related issues
$wsOverzicht.Range("E58").Value = "476, 480, 481, 479"

# Restore the view/selection state recorded for this sheet
$wsOverzicht.Application.ActiveWindow.ScrollRow = 19
$wsOverzicht.Range("E25").Select()

# ---------------------------------------------------------------------------
# Sheet "MVP issues op github (dd 28-05)": mark several github issues with
# whether they are (now) captured in the MVP overview
# ---------------------------------------------------------------------------
$wsGithub = $wb.Worksheets.Item("MVP issues op github (dd 28-05)")

$wsGithub.Range("A96").Value = "Nee (regel 35?)"
$wsGithub.Range("A99").Value = "Ja"
$wsGithub.Range("A102").Value = "Nee (nu wel regel 30)"
$wsGithub.Range("A117").Value = "Ja"
$wsGithub.Range("A120").Value = "Nee (nu wel regel 58, voldoende?"
$wsGithub.Range("A123").Value = "Ja"
$wsGithub.Range("A129").Value = "Nee"
$wsGithub.Range("A132").Value = "Ja"

# Restore the view/selection state recorded for this sheet
$wsGithub.Application.ActiveWindow.ScrollRow = 121
$wsGithub.Range("A135:A136").Select()

$wb.Save()
